$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.648.36'
$ws.Range("E2").Value = '  +1.65%  '

$ws.Range("D3").Value = '2.479.33'
$ws.Range("E3").Value = '  +1.54%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.55'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.00'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.19%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("E8").Value = '  +0.72%  '

$ws.Range("D9").Value = '2.474.46'
$ws.Range("E9").Value = '  +1.19%  '

$ws.Range("E10").Value = '  -0.05%  '

$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.27'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.12%  '

$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.358'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.18'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.15%  '

$ws.Range("E15").Value = '  -2.54%  '

$ws.Range("D16").Value = '2.938.29'
$ws.Range("E16").Value = '  +3.68%  '

$ws.Range("D17").Value = '63.481.78'
$ws.Range("E17").Value = '  +1.71%  '

$ws.Range("D18").Value = '2.485.63'
$ws.Range("E18").Value = '  +1.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.54'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.36'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +5.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '329.77'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.20'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.08'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +19.38%  '

$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.11'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '626.62'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +12.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000104'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +4.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.59'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.95%  '

$ws.Range("D29").Value = '2.619.27'
$ws.Range("E29").Value = '  +2.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.53'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +5.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.38'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.27%  '

$ws.Range("E33").Value = '  -3.37%  '

$ws.Range("E34").Value = '  +1.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.24'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +7.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.52'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.384'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.47'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.81'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '147.81'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.82'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.71'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +12.38%  '

$ws.Range("E44").Value = '  -0.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '149.64'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.76'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.28'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0543'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.605'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0235'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0919'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.15%  '
